$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Loss-of-sale records received for 22-12-2025 through 25-12-2025 are appended
# below the existing table (rows 73-99), continuing the "#" sequence and
# reusing the same per-column layout/formatting as the rows above them.

# Pre-format the Date/Function Date columns as Text for the new rows so that
# ambiguous dd-mm-yyyy strings (e.g. "12-01-2026") are not auto-parsed into
# date serial numbers -- they must stay literal strings, like the existing rows.
$ws.Range("B73:B99").NumberFormat = "@"
$ws.Range("E73:E99").NumberFormat = "@"

# Row 73 (# = 71)
$ws.Cells.Item(73, 1).Value = 71
$ws.Cells.Item(73, 1).NumberFormat = "0"
$ws.Cells.Item(73, 2).Value = '22-12-2025'
$ws.Cells.Item(73, 3).Value = 'Abishek'
$ws.Cells.Item(73, 4).Value = 9847281864
$ws.Cells.Item(73, 4).NumberFormat = "0"
$ws.Cells.Item(73, 5).Value = '12-01-2026'
$ws.Cells.Item(73, 6).Value = 'RAYAN K B'
$ws.Cells.Item(73, 7).Value = 'Loss'
$ws.Cells.Item(73, 8).Value = 'PRODUCT'
$ws.Cells.Item(73, 9).Value = 'REQUIRED MODEL NOT AVAILABLE'
$ws.Cells.Item(73, 10).Value = '-'
$ws.Cells.Item(73, 11).Value = 'PRODUCT ALLREADY BOOKED'

# Row 74 (# = 72)
$ws.Cells.Item(74, 1).Value = 72
$ws.Cells.Item(74, 1).NumberFormat = "0"
$ws.Cells.Item(74, 2).Value = '22-12-2025'
$ws.Cells.Item(74, 3).Value = 'Himal'
$ws.Cells.Item(74, 4).Value = 8547870989
$ws.Cells.Item(74, 4).NumberFormat = "0"
$ws.Cells.Item(74, 5).Value = '26-12-2025'
$ws.Cells.Item(74, 6).Value = 'RAYAN K B'
$ws.Cells.Item(74, 7).Value = 'Loss'
$ws.Cells.Item(74, 8).Value = 'PRODUCT'
$ws.Cells.Item(74, 9).Value = 'REQUIRED MODEL NOT AVAILABLE'
$ws.Cells.Item(74, 10).Value = '-'
$ws.Cells.Item(74, 11).Value = 'REQUIRED MODEL ALLREADY  BOOKED'

# Row 75 (# = 73)
$ws.Cells.Item(75, 1).Value = 73
$ws.Cells.Item(75, 1).NumberFormat = "0"
$ws.Cells.Item(75, 2).Value = '22-12-2025'
$ws.Cells.Item(75, 3).Value = 'SHAJU'
$ws.Cells.Item(75, 4).Value = 8891597962
$ws.Cells.Item(75, 4).NumberFormat = "0"
$ws.Cells.Item(75, 5).Value = '28-12-2025'
$ws.Cells.Item(75, 6).Value = 'ATHULKIRSHNA CS'
$ws.Cells.Item(75, 7).Value = 'Loss'
$ws.Cells.Item(75, 8).Value = 'ENQUIRY'
$ws.Cells.Item(75, 9).Value = 'Enquiry for Relative/Friend'
$ws.Cells.Item(75, 10).Value = '-'
$ws.Cells.Item(75, 11).Value = 'HE FIRST VISITED OUR OUTLET HE NEEDS TO VISIT OTHER SUIT STORES AND HE WILL CONFIRM TOMORROW MORNING.'

# Row 76 (# = 74)
$ws.Cells.Item(76, 1).Value = 74
$ws.Cells.Item(76, 1).NumberFormat = "0"
$ws.Cells.Item(76, 2).Value = '22-12-2025'
$ws.Cells.Item(76, 3).Value = 'AKHIN'
$ws.Cells.Item(76, 4).Value = 9567228804
$ws.Cells.Item(76, 4).NumberFormat = "0"
$ws.Cells.Item(76, 5).Value = '29-12-2025'
$ws.Cells.Item(76, 6).Value = 'RASEEB E A'
$ws.Cells.Item(76, 7).Value = 'Loss'
$ws.Cells.Item(76, 8).Value = 'PRODUCT'
$ws.Cells.Item(76, 9).Value = 'PRODUCT NOT AVAILABLE'
$ws.Cells.Item(76, 10).Value = '-'
$ws.Cells.Item(76, 11).Value = 'OPEN BANGALA BLACK NOT AVAILABLE FOR HIS DATE.'

# Row 77 (# = 75)
$ws.Cells.Item(77, 1).Value = 75
$ws.Cells.Item(77, 1).NumberFormat = "0"
$ws.Cells.Item(77, 2).Value = '22-12-2025'
$ws.Cells.Item(77, 3).Value = 'EMANUAL'
$ws.Cells.Item(77, 4).Value = 9778258992
$ws.Cells.Item(77, 4).NumberFormat = "0"
$ws.Cells.Item(77, 5).Value = '27-12-2025'
$ws.Cells.Item(77, 6).Value = 'RASEEB E A'
$ws.Cells.Item(77, 7).Value = 'Loss'
$ws.Cells.Item(77, 8).Value = 'SIZE NOT SUITABLE'
$ws.Cells.Item(77, 9).Value = 'SIZE TOO SMALL'
$ws.Cells.Item(77, 10).Value = '-'
$ws.Cells.Item(77, 11).Value = 'HE NEED 32 SIZE BLAZER'

# Row 78 (# = 76)
$ws.Cells.Item(78, 1).Value = 76
$ws.Cells.Item(78, 1).NumberFormat = "0"
$ws.Cells.Item(78, 2).Value = '22-12-2025'
$ws.Cells.Item(78, 3).Value = 'shahabas'
$ws.Cells.Item(78, 4).Value = 9037427363
$ws.Cells.Item(78, 4).NumberFormat = "0"
$ws.Cells.Item(78, 5).Value = '24-01-2026'
$ws.Cells.Item(78, 6).Value = 'RASEEB E A'
$ws.Cells.Item(78, 7).Value = 'Loss'
$ws.Cells.Item(78, 8).Value = 'ENQUIRY'
$ws.Cells.Item(78, 9).Value = 'Enquiry for Relative/Friend'
$ws.Cells.Item(78, 10).Value = '-'
$ws.Cells.Item(78, 11).Value = 'JUST VISIT ENQUIRY FOR FRIEND'

# Row 79 (# = 77)
$ws.Cells.Item(79, 1).Value = 77
$ws.Cells.Item(79, 1).NumberFormat = "0"
$ws.Cells.Item(79, 2).Value = '22-12-2025'
$ws.Cells.Item(79, 3).Value = 'Shamsil'
$ws.Cells.Item(79, 4).Value = 9946536516
$ws.Cells.Item(79, 4).NumberFormat = "0"
$ws.Cells.Item(79, 5).Value = '27-12-2025'
$ws.Cells.Item(79, 6).Value = 'RAYAN K B'
$ws.Cells.Item(79, 7).Value = 'Loss'
$ws.Cells.Item(79, 8).Value = 'PRODUCT'
$ws.Cells.Item(79, 9).Value = 'PRODUCT NOT AVAILABLE'
$ws.Cells.Item(79, 10).Value = '-'
$ws.Cells.Item(79, 11).Value = 'PRUDUCT ALLREADY BOOKED'

# Row 80 (# = 78)
$ws.Cells.Item(80, 1).Value = 78
$ws.Cells.Item(80, 1).NumberFormat = "0"
$ws.Cells.Item(80, 2).Value = '23-12-2025'
$ws.Cells.Item(80, 3).Value = 'rinto'
$ws.Cells.Item(80, 4).Value = 8921817705
$ws.Cells.Item(80, 4).NumberFormat = "0"
$ws.Cells.Item(80, 5).Value = '29-12-2025'
$ws.Cells.Item(80, 6).Value = 'RAYAN K B'
$ws.Cells.Item(80, 7).Value = 'Loss'
$ws.Cells.Item(80, 8).Value = 'ENQUIRY'
$ws.Cells.Item(80, 9).Value = 'ENQUIRY WITHOUT BRIDE/FAMILY'
$ws.Cells.Item(80, 10).Value = '-'
$ws.Cells.Item(80, 11).Value = 'UPDATE WITH IN 4 DAYS'

# Row 81 (# = 79)
$ws.Cells.Item(81, 1).Value = 79
$ws.Cells.Item(81, 1).NumberFormat = "0"
$ws.Cells.Item(81, 2).Value = '23-12-2025'
$ws.Cells.Item(81, 3).Value = 'Arun'
$ws.Cells.Item(81, 4).Value = 8129536664
$ws.Cells.Item(81, 4).NumberFormat = "0"
$ws.Cells.Item(81, 5).Value = '17-01-2026'
$ws.Cells.Item(81, 6).Value = 'RAYAN K B'
$ws.Cells.Item(81, 7).Value = 'Loss'
$ws.Cells.Item(81, 8).Value = 'ENQUIRY'
$ws.Cells.Item(81, 9).Value = 'ENQUIRY WITHOUT BRIDE/FAMILY'
$ws.Cells.Item(81, 10).Value = '-'
$ws.Cells.Item(81, 11).Value = 'They liked the product and customer will come tomorrow with bride and family'

# Row 82 (# = 80)
$ws.Cells.Item(82, 1).Value = 80
$ws.Cells.Item(82, 1).NumberFormat = "0"
$ws.Cells.Item(82, 2).Value = '23-12-2025'
$ws.Cells.Item(82, 3).Value = 'NIHAS'
$ws.Cells.Item(82, 4).Value = 9995974196
$ws.Cells.Item(82, 4).NumberFormat = "0"
$ws.Cells.Item(82, 5).Value = '04-01-2026'
$ws.Cells.Item(82, 6).Value = 'SHYAMNADH T J'
$ws.Cells.Item(82, 7).Value = 'Loss'
$ws.Cells.Item(82, 8).Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Cells.Item(82, 9).Value = 'FAMILY DISAPPROVEL'
$ws.Cells.Item(82, 10).Value = '-'
$ws.Cells.Item(82, 11).Value = 'N/J cream 40,42'

# Row 83 (# = 81)
$ws.Cells.Item(83, 1).Value = 81
$ws.Cells.Item(83, 1).NumberFormat = "0"
$ws.Cells.Item(83, 2).Value = '23-12-2025'
$ws.Cells.Item(83, 3).Value = 'lince'
$ws.Cells.Item(83, 4).Value = 9074571235
$ws.Cells.Item(83, 4).NumberFormat = "0"
$ws.Cells.Item(83, 5).Value = '23-12-2025'
$ws.Cells.Item(83, 6).Value = 'ATHULKIRSHNA CS'
$ws.Cells.Item(83, 7).Value = 'Loss'
$ws.Cells.Item(83, 8).Value = 'ENQUIRY'
$ws.Cells.Item(83, 9).Value = 'ENQUIRY WITHOUT BRIDE/FAMILY'
$ws.Cells.Item(83, 10).Value = '-'
$ws.Cells.Item(83, 11).Value = 'customer liked the product and customer will visit with family and book later'

# Row 84 (# = 82)
$ws.Cells.Item(84, 1).Value = 82
$ws.Cells.Item(84, 1).NumberFormat = "0"
$ws.Cells.Item(84, 2).Value = '23-12-2025'
$ws.Cells.Item(84, 3).Value = 'Bibin'
$ws.Cells.Item(84, 4).Value = 9497326932
$ws.Cells.Item(84, 4).NumberFormat = "0"
$ws.Cells.Item(84, 5).Value = '27-12-2025'
$ws.Cells.Item(84, 6).Value = 'ATHULKIRSHNA CS'
$ws.Cells.Item(84, 7).Value = 'Loss'
$ws.Cells.Item(84, 8).Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Cells.Item(84, 9).Value = 'FAMILY DISAPPROVEL'
$ws.Cells.Item(84, 10).Value = '-'
$ws.Cells.Item(84, 11).Value = 'come tomorrow'

# Row 85 (# = 83)
$ws.Cells.Item(85, 1).Value = 83
$ws.Cells.Item(85, 1).NumberFormat = "0"
$ws.Cells.Item(85, 2).Value = '23-12-2025'
$ws.Cells.Item(85, 3).Value = 'Akshay'
$ws.Cells.Item(85, 4).Value = 9048577745
$ws.Cells.Item(85, 4).NumberFormat = "0"
$ws.Cells.Item(85, 5).Value = '28-12-2025'
$ws.Cells.Item(85, 6).Value = 'ATHULKIRSHNA CS'
$ws.Cells.Item(85, 7).Value = 'Loss'
$ws.Cells.Item(85, 8).Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Cells.Item(85, 9).Value = 'FAMILY DISAPPROVEL'
$ws.Cells.Item(85, 10).Value = '-'
$ws.Cells.Item(85, 11).Value = 'visit tomorrow'

# Row 86 (# = 84)
$ws.Cells.Item(86, 1).Value = 84
$ws.Cells.Item(86, 1).NumberFormat = "0"
$ws.Cells.Item(86, 2).Value = '24-12-2025'
$ws.Cells.Item(86, 3).Value = 'JOWIN'
$ws.Cells.Item(86, 4).Value = 7907473358
$ws.Cells.Item(86, 4).NumberFormat = "0"
$ws.Cells.Item(86, 5).Value = '31-12-2025'
$ws.Cells.Item(86, 6).Value = 'SHYAMNADH T J'
$ws.Cells.Item(86, 7).Value = 'Loss'
$ws.Cells.Item(86, 8).Value = 'PRODUCT'
$ws.Cells.Item(86, 9).Value = 'PRODUCT NOT AVAILABLE'
$ws.Cells.Item(86, 10).Value = '-'
$ws.Cells.Item(86, 11).Value = 'LAVENDER COLOUR SIUT'

# Row 87 (# = 85)
$ws.Cells.Item(87, 1).Value = 85
$ws.Cells.Item(87, 1).NumberFormat = "0"
$ws.Cells.Item(87, 2).Value = '24-12-2025'
$ws.Cells.Item(87, 3).Value = 'LEO'
$ws.Cells.Item(87, 4).Value = 8714441119
$ws.Cells.Item(87, 4).NumberFormat = "0"
$ws.Cells.Item(87, 5).Value = '04-01-2026'
$ws.Cells.Item(87, 6).Value = 'SHYAMNADH T J'
$ws.Cells.Item(87, 7).Value = 'Loss'
$ws.Cells.Item(87, 8).Value = 'ENQUIRY'
$ws.Cells.Item(87, 9).Value = 'ENQUIRY WITHOUT BRIDE/FAMILY'
$ws.Cells.Item(87, 10).Value = '-'
$ws.Cells.Item(87, 11).Value = 'HE NEED TO DISCUSS WITH HIS BRIDE'

# Row 88 (# = 86)
$ws.Cells.Item(88, 1).Value = 86
$ws.Cells.Item(88, 1).NumberFormat = "0"
$ws.Cells.Item(88, 2).Value = '24-12-2025'
$ws.Cells.Item(88, 3).Value = 'surave'
$ws.Cells.Item(88, 4).Value = 9048983957
$ws.Cells.Item(88, 4).NumberFormat = "0"
$ws.Cells.Item(88, 5).Value = '29-12-2025'
$ws.Cells.Item(88, 6).Value = 'RASEEB E A'
$ws.Cells.Item(88, 7).Value = 'Loss'
$ws.Cells.Item(88, 8).Value = 'ENQUIRY'
$ws.Cells.Item(88, 9).Value = 'ENQUIRY WITHOUT BRIDE/FAMILY'
$ws.Cells.Item(88, 10).Value = '-'
$ws.Cells.Item(88, 11).Value = 'enquiry Just visit'

# Row 89 (# = 87)
$ws.Cells.Item(89, 1).Value = 87
$ws.Cells.Item(89, 1).NumberFormat = "0"
$ws.Cells.Item(89, 2).Value = '24-12-2025'
$ws.Cells.Item(89, 3).Value = 'ben'
$ws.Cells.Item(89, 4).Value = 8594035331
$ws.Cells.Item(89, 4).NumberFormat = "0"
$ws.Cells.Item(89, 5).Value = '07-01-2026'
$ws.Cells.Item(89, 6).Value = 'RASEEB E A'
$ws.Cells.Item(89, 7).Value = 'Loss'
$ws.Cells.Item(89, 8).Value = 'ENQUIRY'
$ws.Cells.Item(89, 9).Value = 'ENQUIRY WITHOUT TRIAL'
$ws.Cells.Item(89, 10).Value = '-'
$ws.Cells.Item(89, 11).Value = 'jest visit'

# Row 90 (# = 88)
$ws.Cells.Item(90, 1).Value = 88
$ws.Cells.Item(90, 1).NumberFormat = "0"
$ws.Cells.Item(90, 2).Value = '24-12-2025'
$ws.Cells.Item(90, 3).Value = 'BINU'
$ws.Cells.Item(90, 4).Value = 9446278195
$ws.Cells.Item(90, 4).NumberFormat = "0"
$ws.Cells.Item(90, 5).Value = '04-01-2026'
$ws.Cells.Item(90, 6).Value = 'MUHAMMED JASIR. V'
$ws.Cells.Item(90, 7).Value = 'Loss'
$ws.Cells.Item(90, 8).Value = 'ENQUIRY'
$ws.Cells.Item(90, 9).Value = 'ENQUIRY WITHOUT TRIAL'
$ws.Cells.Item(90, 10).Value = '-'
$ws.Cells.Item(90, 11).Value = 'HE COMES WITH FAMILY TRAILED BLACK TEXTURE 40 SIZE BLAZER AND LEAVE.'

# Row 91 (# = 89)
$ws.Cells.Item(91, 1).Value = 89
$ws.Cells.Item(91, 1).NumberFormat = "0"
$ws.Cells.Item(91, 2).Value = '24-12-2025'
$ws.Cells.Item(91, 3).Value = 'Alfred'
$ws.Cells.Item(91, 4).Value = 9656573522
$ws.Cells.Item(91, 4).NumberFormat = "0"
$ws.Cells.Item(91, 5).Value = '07-01-2026'
$ws.Cells.Item(91, 6).Value = 'RAYAN K B'
$ws.Cells.Item(91, 7).Value = 'Loss'
$ws.Cells.Item(91, 8).Value = 'ENQUIRY'
$ws.Cells.Item(91, 9).Value = 'ENQUIRY WITHOUT BRIDE/FAMILY'
$ws.Cells.Item(91, 10).Value = '-'
$ws.Cells.Item(91, 11).Value = 'HE NEEDS TO DISCUSS WITH FAMILY'

# Row 92 (# = 90)
$ws.Cells.Item(92, 1).Value = 90
$ws.Cells.Item(92, 1).NumberFormat = "0"
$ws.Cells.Item(92, 2).Value = '24-12-2025'
$ws.Cells.Item(92, 3).Value = 'JOHN'
$ws.Cells.Item(92, 4).Value = 7736692808
$ws.Cells.Item(92, 4).NumberFormat = "0"
$ws.Cells.Item(92, 5).Value = '03-01-2026'
$ws.Cells.Item(92, 6).Value = 'ATHULKIRSHNA CS'
$ws.Cells.Item(92, 7).Value = 'Loss'
$ws.Cells.Item(92, 8).Value = 'ENQUIRY'
$ws.Cells.Item(92, 9).Value = 'ENQUIRY WITHOUT TRIAL'
$ws.Cells.Item(92, 10).Value = '-'
$ws.Cells.Item(92, 11).Value = 'JUST VISIT FOR LOOKING RENTAL'

# Row 93 (# = 91)
$ws.Cells.Item(93, 1).Value = 91
$ws.Cells.Item(93, 1).NumberFormat = "0"
$ws.Cells.Item(93, 2).Value = '24-12-2025'
$ws.Cells.Item(93, 3).Value = 'Gowtham'
$ws.Cells.Item(93, 4).Value = 8943210145
$ws.Cells.Item(93, 4).NumberFormat = "0"
$ws.Cells.Item(93, 5).Value = '29-12-2025'
$ws.Cells.Item(93, 6).Value = 'RASEEB E A'
$ws.Cells.Item(93, 7).Value = 'Loss'
$ws.Cells.Item(93, 8).Value = 'SIZE NOT SUITABLE'
$ws.Cells.Item(93, 9).Value = 'SIZE TOO LARGE'
$ws.Cells.Item(93, 10).Value = '-'
$ws.Cells.Item(93, 11).Value = '46 SIZE KURTHA NEEDED'

# Row 94 (# = 92)
$ws.Cells.Item(94, 1).Value = 92
$ws.Cells.Item(94, 1).NumberFormat = "0"
$ws.Cells.Item(94, 2).Value = '25-12-2025'
$ws.Cells.Item(94, 3).Value = 'sidheek'
$ws.Cells.Item(94, 4).Value = 8129199845
$ws.Cells.Item(94, 4).NumberFormat = "0"
$ws.Cells.Item(94, 5).Value = '28-12-2025'
$ws.Cells.Item(94, 6).Value = 'RASEEB E A'
$ws.Cells.Item(94, 7).Value = 'Loss'
$ws.Cells.Item(94, 8).Value = 'ENQUIRY'
$ws.Cells.Item(94, 9).Value = 'ENQUIRY WITHOUT TRIAL'
$ws.Cells.Item(94, 10).Value = '-'
$ws.Cells.Item(94, 11).Value = 'just visit'

# Row 95 (# = 93)
$ws.Cells.Item(95, 1).Value = 93
$ws.Cells.Item(95, 1).NumberFormat = "0"
$ws.Cells.Item(95, 2).Value = '25-12-2025'
$ws.Cells.Item(95, 3).Value = 'Suman'
$ws.Cells.Item(95, 4).Value = 9567059770
$ws.Cells.Item(95, 4).NumberFormat = "0"
$ws.Cells.Item(95, 5).Value = '28-12-2025'
$ws.Cells.Item(95, 6).Value = 'ATHULKIRSHNA CS'
$ws.Cells.Item(95, 7).Value = 'Loss'
$ws.Cells.Item(95, 8).Value = 'PRODUCT'
$ws.Cells.Item(95, 9).Value = 'Product Already Booked'
$ws.Cells.Item(95, 10).Value = '-'
$ws.Cells.Item(95, 11).Value = 'product not available'

# Row 96 (# = 94)
$ws.Cells.Item(96, 1).Value = 94
$ws.Cells.Item(96, 1).NumberFormat = "0"
$ws.Cells.Item(96, 2).Value = '25-12-2025'
$ws.Cells.Item(96, 3).Value = 'jesleo'
$ws.Cells.Item(96, 4).Value = 9496751850
$ws.Cells.Item(96, 4).NumberFormat = "0"
$ws.Cells.Item(96, 5).Value = '29-12-2025'
$ws.Cells.Item(96, 6).Value = 'ATHULKIRSHNA CS'
$ws.Cells.Item(96, 7).Value = 'Loss'
$ws.Cells.Item(96, 8).Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Cells.Item(96, 9).Value = 'FAMILY DISAPPROVEL'
$ws.Cells.Item(96, 10).Value = '-'
$ws.Cells.Item(96, 11).Value = 'update later'

# Row 97 (# = 95)
$ws.Cells.Item(97, 1).Value = 95
$ws.Cells.Item(97, 1).NumberFormat = "0"
$ws.Cells.Item(97, 2).Value = '25-12-2025'
$ws.Cells.Item(97, 3).Value = 'Munas'
$ws.Cells.Item(97, 4).Value = 9645906666
$ws.Cells.Item(97, 4).NumberFormat = "0"
$ws.Cells.Item(97, 5).Value = '29-12-2025'
$ws.Cells.Item(97, 6).Value = 'RAYAN K B'
$ws.Cells.Item(97, 7).Value = 'Loss'
$ws.Cells.Item(97, 8).Value = 'ENQUIRY'
$ws.Cells.Item(97, 9).Value = 'Enquiry for Relative/Friend'
$ws.Cells.Item(97, 10).Value = '-'
$ws.Cells.Item(97, 11).Value = 'HE WILL COME WITH HIS COUSINS TOMORROW OR NEXT DAY.'

# Row 98 (# = 96)
$ws.Cells.Item(98, 1).Value = 96
$ws.Cells.Item(98, 1).NumberFormat = "0"
$ws.Cells.Item(98, 2).Value = '25-12-2025'
$ws.Cells.Item(98, 3).Value = 'haize'
$ws.Cells.Item(98, 4).Value = 9995125150
$ws.Cells.Item(98, 4).NumberFormat = "0"
$ws.Cells.Item(98, 5).Value = '15-01-2026'
$ws.Cells.Item(98, 6).Value = 'ATHULKIRSHNA CS'
$ws.Cells.Item(98, 7).Value = 'Loss'
$ws.Cells.Item(98, 8).Value = 'ENQUIRY'
$ws.Cells.Item(98, 9).Value = 'ENQUIRY WITHOUT TRIAL'
$ws.Cells.Item(98, 10).Value = '-'
$ws.Cells.Item(98, 11).Value = 'just a enquiry'

# Row 99 (# = 97)
$ws.Cells.Item(99, 1).Value = 97
$ws.Cells.Item(99, 1).NumberFormat = "0"
$ws.Cells.Item(99, 2).Value = '25-12-2025'
$ws.Cells.Item(99, 3).Value = 'AASIL'
$ws.Cells.Item(99, 4).Value = 9947535233
$ws.Cells.Item(99, 4).NumberFormat = "0"
$ws.Cells.Item(99, 5).Value = '11-01-2026'
$ws.Cells.Item(99, 6).Value = 'SHYAMNADH T J'
$ws.Cells.Item(99, 7).Value = 'Loss'
$ws.Cells.Item(99, 8).Value = 'CUSTOMER INTERNAL ISSUES'
$ws.Cells.Item(99, 9).Value = 'BUDGET RESTRICTIONS'
$ws.Cells.Item(99, 10).Value = '-'
$ws.Cells.Item(99, 11).Value = 'THEY CHOOSED TOTAL OF 9 QTY VALUE ALMOST 35K AND THEY NEED ADDITIONAL 7 PRODUCTS IT ALMOST HITS 65-70K RANGE BUT THE GROOM IS THE SPONSER HIS BUDGET IS MAXIMUM 40K RANGE.'

# Restore General formatting on the Date/Function Date columns now that the
# literal text has been committed, so the cells carry no explicit style -
# matching the plain inlineStr cells used throughout the rest of the sheet.
$ws.Range("B73:B99").Style = "Normal"
$ws.Range("E73:E99").Style = "Normal"

# Column F ("Staff") auto-resizes because the new row for BINU includes the
# longest staff name seen so far, "MUHAMMED JASIR. V" (17 characters), which
# is wider than the previous longest entry, "ATHULKIRSHNA CS" (15 characters).
$ws.Columns.Item(6).ColumnWidth = 22.09
